# Generate Report for Handback
#
# This script brings the localization-status workbook from "handoff pending"
# to "handed back" state:
#   1. Every cell showing the old "Ready for handoff" status text is updated
#      to "Handed back: in sync with en-US" (Overview + zh-cn + de-de sheets).
#   2. The (previously blank "0001-01-01 00:00:00") Latest Handback DateTime
#      cells get real handback timestamps - different per language sheet.
#   3. The previously-empty "Latest Target File" / "Latest Handback File"
#      columns (F/G) are populated with the handed-back file names, each
#      linked back to the same targets as the corresponding handoff-file
#      hyperlinks already on the row (columns A / D).

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
# ---------------------------------------------------------------------

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusNew
$wsOverview.Range("C2").Value = $statusNew
$wsOverview.Range("B3").Value = $statusNew
$wsOverview.Range("C3").Value = $statusNew

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusNew
$wsZhCn.Range("C3").Value = $statusNew

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusNew
$wsDeDe.Range("C3").Value = $statusNew

# ---------------------------------------------------------------------
# 2. zh-cn sheet: Latest Target File / Latest Handback File / Latest
#    Handback DateTime for both data rows.
# ---------------------------------------------------------------------

$zhCnMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/aa8df0bc51d78e2ce01e8daa87be91cdfea6f30c/e2e/5e50956d-4f7f-49c4-890a-14a77f2b3ae5.md"
$zhCnXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7def112d17ac5f1d29362718ded458a32b25976f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/5e50956d-4f7f-49c4-890a-14a77f2b3ae5.def9e0405441e242ba97c006128289193b32ba53.zh-cn.xlf"
$zhCnMdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/aa8df0bc51d78e2ce01e8daa87be91cdfea6f30c/e2e/829de9ef-8911-4647-94cf-9c56b33d2d13.md"
$zhCnXlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7def112d17ac5f1d29362718ded458a32b25976f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/829de9ef-8911-4647-94cf-9c56b33d2d13.b2be44d7883d4a4e516ce06e34a2d751d89c0d22.zh-cn.xlf"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), $zhCnMdUrl, "", "", "5e50956d-4f7f-49c4-890a-14a77f2b3ae5.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), $zhCnXlfUrl, "", "", "5e50956d-4f7f-49c4-890a-14a77f2b3ae5.def9e0405441e242ba97c006128289193b32ba53.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), $zhCnMdUrl2, "", "", "829de9ef-8911-4647-94cf-9c56b33d2d13.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), $zhCnXlfUrl2, "", "", "829de9ef-8911-4647-94cf-9c56b33d2d13.b2be44d7883d4a4e516ce06e34a2d751d89c0d22.zh-cn.xlf")

$wsZhCn.Range("H2").Value = "2016-03-19 00:14:31"
$wsZhCn.Range("H3").Value = "2016-03-19 00:14:31"

# ---------------------------------------------------------------------
# 3. de-de sheet: Latest Target File / Latest Handback File / Latest
#    Handback DateTime for both data rows.
# ---------------------------------------------------------------------

$deDeMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/aa8df0bc51d78e2ce01e8daa87be91cdfea6f30c/e2e/5e50956d-4f7f-49c4-890a-14a77f2b3ae5.md"
$deDeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/450a4302800d4564497bae8e7bd4fa41ae05eda3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/5e50956d-4f7f-49c4-890a-14a77f2b3ae5.def9e0405441e242ba97c006128289193b32ba53.de-de.xlf"
$deDeMdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/aa8df0bc51d78e2ce01e8daa87be91cdfea6f30c/e2e/829de9ef-8911-4647-94cf-9c56b33d2d13.md"
$deDeXlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/450a4302800d4564497bae8e7bd4fa41ae05eda3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/829de9ef-8911-4647-94cf-9c56b33d2d13.b2be44d7883d4a4e516ce06e34a2d751d89c0d22.de-de.xlf"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), $deDeMdUrl, "", "", "5e50956d-4f7f-49c4-890a-14a77f2b3ae5.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), $deDeXlfUrl, "", "", "5e50956d-4f7f-49c4-890a-14a77f2b3ae5.def9e0405441e242ba97c006128289193b32ba53.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), $deDeMdUrl2, "", "", "829de9ef-8911-4647-94cf-9c56b33d2d13.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), $deDeXlfUrl2, "", "", "829de9ef-8911-4647-94cf-9c56b33d2d13.b2be44d7883d4a4e516ce06e34a2d751d89c0d22.de-de.xlf")

$wsDeDe.Range("H2").Value = "2016-03-19 00:14:35"
$wsDeDe.Range("H3").Value = "2016-03-19 00:14:35"
